$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Remaining Tasks")

# 1) Task #4 ("Adding Products") now has a working example: update Status + Note
$ws.Range("D8").Value = "Has Working Example"
$ws.Range("E8").Value = "Have to read form as JSON data"

# 2) Preserve the existing review comment's text before shifting rows
$commentText = $ws.Range("D21").Comment.Text()
$ws.Range("D21").Comment.Delete()

# 3) Remove the now-obsolete "Putting Products On Promotion" task row (old row 10)
$ws.Rows("10:10").Delete()

# 4) Renumber the "Number" column for every task row that shifted up
for ($r = 10; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 4
}

# 5) Re-attach the comment to its new location (one row up, D20)
$ws.Range("D20").AddComment($commentText)

# 5) Refresh the view (matches what Excel records after this edit/scroll)
$ws.Range("G7").Select()
